$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text
$ws.Range("A1").Value = "Peak Loads for March 2025 (Mthembanji)"

# Update values for days 1-3
$ws.Range("B3").Value = 2.0015744166
$ws.Range("B4").Value = 2.0794223832
$ws.Range("B5").Value = 0.709

# Delete rows 6 through 12 (days 4-10), shifting cells up
$ws.Range("A6:C12").EntireRow.Delete()
